$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy B1's formatting (style) onto C1, then set its value to "Idade"
$ws.Range("B1").Copy($ws.Range("C1"))
$ws.Range("C1").Value = "Idade"

# Delete row 2 entirely (removes "as" / "asssss" entries)
$ws.Rows.Item(2).Delete()
